$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 11495
$ws.Range("I51").Value = 11495
$ws.Range("K51").Value = 11495
$ws.Range("M51").Value = -11011

$ws.Range("H69").Value = 7020.75
$ws.Range("J69").Value = 7020.75
$ws.Range("L69").Value = 21062.25
$ws.Range("N69").Value = -22810.25

$ws.Range("H72").Value = 7020.75
$ws.Range("J72").Value = 7020.75
$ws.Range("L72").Value = 63186.75
$ws.Range("N72").Value = -71922.75

$ws.Range("H80").Value = 5682399.5
$ws.Range("I80").Value = 10417001
$ws.Range("J80").Value = 877.2
$ws.Range("K80").Value = 31251003
$ws.Range("L80").Value = 2631.6
$ws.Range("M80").Value = -31250005
$ws.Range("N80").Value = -4627.6

$ws.Range("H83").Value = 5682399.5
$ws.Range("I83").Value = 10417001
$ws.Range("J83").Value = 877.2
$ws.Range("K83").Value = 93753009
$ws.Range("L83").Value = 7894.8
$ws.Range("M83").Value = -93748017
$ws.Range("N83").Value = -17878.8

$ws.Range("H100").Value = 5264.381
$ws.Range("I100").Value = 2680.6428
$ws.Range("J100").Value = 10431.857
$ws.Range("K100").Value = 2680.6428
$ws.Range("L100").Value = 10431.857
$ws.Range("M100").Value = -2139.6428
$ws.Range("N100").Value = -11513.857

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 0

$ws.Range("H135").Value = 5872.273
$ws.Range("I135").Value = 2614
$ws.Range("J135").Value = 20534.5
$ws.Range("K135").Value = 23526
$ws.Range("L135").Value = 184810.5
$ws.Range("M135").Value = -20991
$ws.Range("N135").Value = -189880.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2523.3572
$ws.Range("I45").Value = 2332.8462
$ws.Range("K45").Value = 2332.8462
$ws.Range("M45").Value = -1955.8462

$ws.Range("H122").Value = 1358.4375
$ws.Range("I122").Value = 1303.6154
$ws.Range("J122").Value = 1596
$ws.Range("K122").Value = 3910.8462
$ws.Range("L122").Value = 4788
$ws.Range("M122").Value = -1460.8462
$ws.Range("N122").Value = -9688

$ws.Range("H132").Value = 1980.6428
$ws.Range("I132").Value = 1706.18
$ws.Range("K132").Value = 5118.54
$ws.Range("M132").Value = -2588.54

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1616.75
$ws.Range("I86").Value = 1486.875
$ws.Range("K86").Value = 1486.875
$ws.Range("M86").Value = -363.875

$ws.Range("H89").Value = 1616.75
$ws.Range("I89").Value = 1486.875
$ws.Range("K89").Value = 7434.375
$ws.Range("M89").Value = -1818.375

$ws.Range("H107").Value = 2021.8182
$ws.Range("J107").Value = 1872
$ws.Range("L107").Value = 1872
$ws.Range("N107").Value = -5712

$ws.Range("H134").Value = 2428.7415
$ws.Range("I134").Value = 1725.7097
$ws.Range("K134").Value = 5177.1291
$ws.Range("M134").Value = -2642.1291

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1827.0513
$ws.Range("I58").Value = 1165.5714
$ws.Range("K58").Value = 1165.5714
$ws.Range("M58").Value = -962.5714

$ws.Range("H62").Value = 2589.1
$ws.Range("J62").Value = 2853
$ws.Range("L62").Value = 2853
$ws.Range("N62").Value = -4101

$ws.Range("H65").Value = 2589.1
$ws.Range("J65").Value = 2853
$ws.Range("L65").Value = 14265
$ws.Range("N65").Value = -20505

$ws.Range("H99").Value = 3379.8
$ws.Range("I99").Value = 3000
$ws.Range("J99").Value = 3474.75
$ws.Range("K99").Value = 3000
$ws.Range("L99").Value = 3474.75
$ws.Range("M99").Value = -1502
$ws.Range("N99").Value = -6470.75

$ws.Range("H107").Value = 693
$ws.Range("I107").Value = 631.125
$ws.Range("J107").Value = 816.75
$ws.Range("K107").Value = 631.125
$ws.Range("L107").Value = 816.75
$ws.Range("M107").Value = 1288.875
$ws.Range("N107").Value = -4656.75

$ws.Range("H126").Value = 3379.8
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 3474.75
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 10424.25
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -15364.25

$ws.Range("H132").Value = 2501500
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H136").Value = 1827.0513
$ws.Range("I136").Value = 1165.5714
$ws.Range("K136").Value = 3496.7142
$ws.Range("M136").Value = -946.7142000000003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 156293.5
$ws.Range("J2").Value = 87.42856999999999
$ws.Range("L2").Value = 524.57142
$ws.Range("N2").Value = -750.57142

$ws.Range("H38").Value = 10.8
$ws.Range("J38").Value = 13.142858
$ws.Range("L38").Value = 39.428574
$ws.Range("N38").Value = -733.428574

$ws.Range("H75").Value = 4284.1055
$ws.Range("J75").Value = 4825.533
$ws.Range("L75").Value = 14476.599
$ws.Range("N75").Value = -16472.599

$ws.Range("H78").Value = 4284.1055
$ws.Range("J78").Value = 4825.533
$ws.Range("L78").Value = 43429.79700000001
$ws.Range("N78").Value = -53413.79700000001

$ws.Range("H92").Value = 784.2857
$ws.Range("I92").Value = 738
$ws.Range("K92").Value = 2214
$ws.Range("M92").Value = -966

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 267311.94
$ws.Range("I80").Value = 558872.3
$ws.Range("J80").Value = 4907.6
$ws.Range("K80").Value = 558872.3
$ws.Range("L80").Value = 4907.6
$ws.Range("M80").Value = -557874.3
$ws.Range("N80").Value = -6903.6

$ws.Range("H83").Value = 267311.94
$ws.Range("I83").Value = 558872.3
$ws.Range("J83").Value = 4907.6
$ws.Range("K83").Value = 2794361.5
$ws.Range("L83").Value = 24538
$ws.Range("M83").Value = -2789369.5
$ws.Range("N83").Value = -34522

$ws.Range("H107").Value = 1715.6154
$ws.Range("J107").Value = 1912.8334
$ws.Range("L107").Value = 1912.8334
$ws.Range("N107").Value = -5752.8334

$ws.Range("H122").Value = 1998.56
$ws.Range("I122").Value = 1503.762
$ws.Range("J122").Value = 4596.25
$ws.Range("K122").Value = 4511.286
$ws.Range("L122").Value = 13788.75
$ws.Range("M122").Value = -2061.286
$ws.Range("N122").Value = -18688.75

$ws.Range("H123").Value = 59800
$ws.Range("J123").Value = 59800
$ws.Range("L123").Value = 59800
$ws.Range("N123").Value = -64700

$ws.Range("H126").Value = 2659.2727
$ws.Range("I126").Value = 2563.4375
$ws.Range("J126").Value = 2914.8333
$ws.Range("K126").Value = 7690.3125
$ws.Range("L126").Value = 8744.499899999999
$ws.Range("M126").Value = -5220.3125
$ws.Range("N126").Value = -13684.4999

$ws.Range("H132").Value = 28581708
$ws.Range("J132").Value = 25867.334
$ws.Range("L132").Value = 77602.00199999999
$ws.Range("N132").Value = -82662.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1057.4445
$ws.Range("I22").Value = 829.75
$ws.Range("J22").Value = 1239.6
$ws.Range("K22").Value = 829.75
$ws.Range("L22").Value = 1239.6
$ws.Range("M22").Value = -534.75
$ws.Range("N22").Value = -1829.6

$ws.Range("H27").Value = 1057.4445
$ws.Range("I27").Value = 829.75
$ws.Range("J27").Value = 1239.6
$ws.Range("K27").Value = 829.75
$ws.Range("L27").Value = 1239.6
$ws.Range("M27").Value = -722.75
$ws.Range("N27").Value = -1453.6

$ws.Range("H46").Value = 1411.1111
$ws.Range("J46").Value = 1616.8334
$ws.Range("L46").Value = 1616.8334
$ws.Range("N46").Value = -1992.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 111392616
$ws.Range("J4").Value = 166671540
$ws.Range("L4").Value = 166671540
$ws.Range("N4").Value = -166671766

$ws.Range("H126").Value = 2114.739
$ws.Range("I126").Value = 2005.8667
$ws.Range("J126").Value = 2318.875
$ws.Range("K126").Value = 6017.6001
$ws.Range("L126").Value = 6956.625
$ws.Range("M126").Value = -3547.6001
$ws.Range("N126").Value = -11896.625

$ws.Range("H132").Value = 1952.0322
$ws.Range("I132").Value = 1427.2693
$ws.Range("K132").Value = 4281.8079
$ws.Range("M132").Value = -1751.8079

$ws.Range("H136").Value = 3920.4866
$ws.Range("I136").Value = 3213.879
$ws.Range("K136").Value = 9641.636999999999
$ws.Range("M136").Value = -7091.636999999999
